$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'307.45"
$ws.Range("E2").Value = "'-2.60%"
$ws.Range("G2").Value = "'13"
$ws.Range("D3").Value = "'37.57"
$ws.Range("E3").Value = "'-5.02%"
$ws.Range("G3").Value = "'13"
$ws.Range("D4").Value = "'5.097"
$ws.Range("E4").Value = "'-0.66%"
$ws.Range("G4").Value = "'13"
$ws.Range("D5").Value = "'0.07878"
$ws.Range("E5").Value = "'-3.95%"
$ws.Range("G5").Value = "'13"
$ws.Range("D6").Value = "'1.983"
$ws.Range("E6").Value = "'1.19%"
$ws.Range("G6").Value = "'13"
$ws.Range("D7").Value = "'4.334"
$ws.Range("E7").Value = "'1.78%"
$ws.Range("G7").Value = "'13"
$ws.Range("D8").Value = "'8.195"
$ws.Range("E8").Value = "'-0.43%"
$ws.Range("G8").Value = "'13"
$ws.Range("D9").Value = "'3.176"
$ws.Range("E9").Value = "'-3.84%"
$ws.Range("G9").Value = "'13"
$ws.Range("D10").Value = "'0.9230"
$ws.Range("E10").Value = "'-0.63%"
$ws.Range("G10").Value = "'13"
$ws.Range("D11").Value = "'0.1289"
$ws.Range("E11").Value = "'-8.57%"
$ws.Range("G11").Value = "'13"
$ws.Range("D12").Value = "'0.1890"
$ws.Range("E12").Value = "'-4.19%"
$ws.Range("G12").Value = "'13"
$ws.Range("D13").Value = "'0.08758"
$ws.Range("E13").Value = "'-3.86%"
$ws.Range("G13").Value = "'13"
$ws.Range("D14").Value = "'0.03424"
$ws.Range("E14").Value = "'-2.39%"
$ws.Range("G14").Value = "'13"
$ws.Range("D15").Value = "'0.09749"
$ws.Range("E15").Value = "'-0.56%"
$ws.Range("G15").Value = "'13"
$ws.Range("E16").Value = "'-0.84%"
$ws.Range("G16").Value = "'13"
$ws.Range("D17").Value = "'0.005917"
$ws.Range("E17").Value = "'-0.03%"
$ws.Range("G17").Value = "'13"
$ws.Range("E18").Value = "'1,782.95%"
$ws.Range("G18").Value = "'13"
$ws.Range("D19").Value = "'3.583"
$ws.Range("E19").Value = "'-1.97%"
$ws.Range("G19").Value = "'13"
$ws.Range("D20").Value = "'0.3437"
$ws.Range("E20").Value = "'-0.74%"
$ws.Range("G20").Value = "'13"
$ws.Range("D21").Value = "'0.1283"
$ws.Range("E21").Value = "'-0.81%"
$ws.Range("G21").Value = "'13"
$ws.Range("D22").Value = "'5.007"
$ws.Range("E22").Value = "'1.74%"
$ws.Range("G22").Value = "'13"
$ws.Range("D23").Value = "'0.2496"
$ws.Range("E23").Value = "'2.24%"
$ws.Range("G23").Value = "'13"
$ws.Range("D24").Value = "'0.04325"
$ws.Range("E24").Value = "'-0.25%"
$ws.Range("G24").Value = "'13"
$ws.Range("D25").Value = "'0.001222"
$ws.Range("E25").Value = "'0.19%"
$ws.Range("G25").Value = "'13"
$ws.Range("D26").Value = "'0.004596"
$ws.Range("E26").Value = "'-4.15%"
$ws.Range("G26").Value = "'13"
$ws.Range("E27").Value = "'177.31%"
$ws.Range("G27").Value = "'13"
$ws.Range("G28").Value = "'13"
$ws.Range("G29").Value = "'13"
$ws.Range("G30").Value = "'13"
$ws.Range("G31").Value = "'13"
$ws.Range("G32").Value = "'13"
$ws.Range("G33").Value = "'13"
$ws.Range("G34").Value = "'13"
$ws.Range("G35").Value = "'13"
$ws.Range("G36").Value = "'13"
$ws.Range("G37").Value = "'13"
$ws.Range("G38").Value = "'13"
$ws.Range("D39").Value = "'0.02306"
$ws.Range("E39").Value = "'2.93%"
$ws.Range("G39").Value = "'13"
$ws.Range("D40").Value = "'0.05022"
$ws.Range("E40").Value = "'-4.75%"
$ws.Range("G40").Value = "'13"
$ws.Range("D41").Value = "'0.007505"
$ws.Range("E41").Value = "'-1.00%"
$ws.Range("G41").Value = "'13"
$ws.Range("D42").Value = "'0.009927"
$ws.Range("E42").Value = "'1.31%"
$ws.Range("G42").Value = "'13"
$ws.Range("D43").Value = "'0.1353"
$ws.Range("E43").Value = "'-1.68%"
$ws.Range("G43").Value = "'13"
$ws.Range("D44").Value = "'0.002093"
$ws.Range("E44").Value = "'-0.49%"
$ws.Range("G44").Value = "'13"
$ws.Range("D45").Value = "'0.008026"
$ws.Range("G45").Value = "'13"
$ws.Range("D46").Value = "'0.00006383"
$ws.Range("E46").Value = "'0.60%"
$ws.Range("G46").Value = "'13"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'0.48%"
$ws.Range("G47").Value = "'13"
$ws.Range("D48").Value = "'0.002997"
$ws.Range("E48").Value = "'8.61%"
$ws.Range("G48").Value = "'13"
$ws.Range("D49").Value = "'0.001202"
$ws.Range("E49").Value = "'0.49%"
$ws.Range("G49").Value = "'13"
$ws.Range("D50").Value = "'0.00002103"
$ws.Range("E50").Value = "'0.48%"
$ws.Range("G50").Value = "'13"
$ws.Range("D51").Value = "'0.0002003"
$ws.Range("E51").Value = "'0.48%"
$ws.Range("G51").Value = "'13"
